$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new data rows 26-35 (new PNG map/chart filenames + their descriptions) ---

# Date column first
$ws.Range("A26:A35").Value = 43900

# Column B (filenames) - new unique strings appended to the shared string table
# in this order (indices 51-60)
$ws.Range("B26").Value = "sprfd_mo_msa_ct_04_spatAuto_global_blk.png"
$ws.Range("B27").Value = "sprfd_mo_msa_ct_04_spatAuto_global_edtot.png"
$ws.Range("B28").Value = "sprfd_mo_msa_ct_04_spatAuto_global_mhi.png"
$ws.Range("B29").Value = "sprfd_mo_msa_ct_04_spatAuto_global_nhi.png"
$ws.Range("B30").Value = "sprfd_mo_msa_ct_04_spatAuto_global_pblk.png"
$ws.Range("B31").Value = "sprfd_mo_msa_ct_04_spatAuto_global_pov.png"
$ws.Range("B32").Value = "sprfd_mo_msa_ct_04_spatAuto_global_pwht.png"
$ws.Range("B33").Value = "sprfd_mo_msa_ct_04_spatAuto_global_wht.png"
$ws.Range("B34").Value = "sprfd_mo_msa_grid02_spatAuto_global_blk.png"
$ws.Range("B35").Value = "sprfd_mo_msa_grid02_spatAuto_global_wht.png"

# Column C (descriptions) - new unique strings appended to the shared string table
# in this order (indices 61-69)
$ws.Range("C26").Value = "Univariate global Moran's I for white population at census tract level."
$ws.Range("C29").Value = "Univariate global Moran's I for no health insurance at census tract level."
$ws.Range("C28").Value = "Univariate global Moran's I for median household income at census tract level."
$ws.Range("C27").Value = "Univariate global Moran's I for total education at census tract level."
$ws.Range("C30").Value = "Univariate global Moran's I for percent black population at census tract level."
$ws.Range("C31").Value = "Univariate global Moran's I for poverty level at census tract level."
$ws.Range("C32").Value = "Univariate global Moran's I for percent white population at census tract level."
$ws.Range("C33").Value = "Univariate global Moran's I for white population at census tract level."
$ws.Range("C34").Value = "Univariate global Moran's I for black population at grid level."
$ws.Range("C35").Value = "Univariate global Moran's I for white population at grid level."

# --- Formatting updates ---
# Column B width widened (~25.7 -> ~50.7 characters) to fit the longer filenames
$ws.Columns.Item(2).ColumnWidth = 49.86

# Row 3 height reduced from 45 to 30
$ws.Rows.Item(3).RowHeight = 30

# --- View update: move active selection down past the new rows ---
$ws.Range("A36").Select() | Out-Null
